$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / percentage updates (never numeric-looking, safe to assign directly)
$ws.Cells.Item(2, 5).Value = "  -0.33%  "
$ws.Cells.Item(3, 5).Value = "  -1.39%  "
$ws.Cells.Item(4, 5).Value = "  +0.03%  "
$ws.Cells.Item(5, 5).Value = "  +0.26%  "
$ws.Cells.Item(6, 5).Value = "  -3.70%  "
$ws.Cells.Item(7, 5).Value = "  +0.12%  "
$ws.Cells.Item(8, 5).Value = "  -1.30%  "
$ws.Cells.Item(9, 5).Value = "  -1.57%  "
$ws.Cells.Item(10, 5).Value = "  -0.70%  "
$ws.Cells.Item(11, 5).Value = "  -2.55%  "
$ws.Cells.Item(12, 5).Value = "  -1.75%  "
$ws.Cells.Item(13, 5).Value = "  -2.97%  "
$ws.Cells.Item(15, 5).Value = "  -2.30%  "
$ws.Cells.Item(16, 5).Value = "  -1.38%  "
$ws.Cells.Item(17, 5).Value = "  -0.46%  "
$ws.Cells.Item(18, 5).Value = "  -1.96%  "
$ws.Cells.Item(19, 5).Value = "  -1.39%  "
$ws.Cells.Item(20, 5).Value = "  +0.92%  "
$ws.Cells.Item(21, 5).Value = "  +1.40%  "
$ws.Cells.Item(22, 5).Value = "  -0.16%  "
$ws.Cells.Item(23, 5).Value = "  -3.75%  "
$ws.Cells.Item(24, 5).Value = "  -1.58%  "
$ws.Cells.Item(26, 5).Value = "  -3.73%  "
$ws.Cells.Item(27, 5).Value = "  +1.76%  "
$ws.Cells.Item(28, 5).Value = "  -0.03%  "
$ws.Cells.Item(29, 5).Value = "  -3.56%  "
$ws.Cells.Item(30, 5).Value = "  -5.24%  "
$ws.Cells.Item(31, 5).Value = "  -1.60%  "
$ws.Cells.Item(32, 5).Value = "  -3.28%  "
$ws.Cells.Item(33, 5).Value = "  -2.25%  "
$ws.Cells.Item(34, 5).Value = "  -7.01%  "
$ws.Cells.Item(35, 5).Value = "  +0.03%  "
$ws.Cells.Item(36, 5).Value = "  -3.25%  "
$ws.Cells.Item(37, 5).Value = "  -2.68%  "
$ws.Cells.Item(38, 5).Value = "  -3.00%  "
$ws.Cells.Item(39, 5).Value = "  +0.65%  "
$ws.Cells.Item(40, 5).Value = "  -5.46%  "
$ws.Cells.Item(41, 5).Value = "  -3.54%  "
$ws.Cells.Item(42, 5).Value = "  -3.75%  "
$ws.Cells.Item(43, 5).Value = "  -1.66%  "
$ws.Cells.Item(44, 5).Value = "  -2.56%  "
$ws.Cells.Item(45, 2).Value = "VeChain"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(45, 5).Value = "  -3.55%  "
$ws.Cells.Item(46, 2).Value = "Monero"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(46, 5).Value = "  -0.50%  "
$ws.Cells.Item(47, 2).Value = "dogwifhat"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(47, 5).Value = "  -2.70%  "
$ws.Cells.Item(49, 5).Value = "  -1.99%  "
$ws.Cells.Item(50, 5).Value = "  -1.63%  "
$ws.Cells.Item(51, 5).Value = "  -3.32%  "

# Price column updates: use a leading apostrophe to force text type (matches
# original inlineStr cells), then ClearFormats to drop the auto-added
# quote-prefix style so the cell style stays at the original default (no "s" attr).
$ws.Cells.Item(2, 4).Value = "'66.802.09"
$ws.Cells.Item(2, 4).ClearFormats()
$ws.Cells.Item(3, 4).Value = "'3.073.89"
$ws.Cells.Item(3, 4).ClearFormats()
$ws.Cells.Item(5, 4).Value = "'577.23"
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(6, 4).Value = "'167.67"
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(8, 4).Value = "'3.073.33"
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(9, 4).Value = "'0.513"
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(10, 4).Value = "'6.38"
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(11, 4).Value = "'0.149"
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(12, 4).Value = "'0.470"
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(14, 4).Value = "'36.05"
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(16, 4).Value = "'3.585.57"
$ws.Cells.Item(16, 4).ClearFormats()
$ws.Cells.Item(17, 4).Value = "'66.752.51"
$ws.Cells.Item(17, 4).ClearFormats()
$ws.Cells.Item(18, 4).Value = "'6.98"
$ws.Cells.Item(18, 4).ClearFormats()
$ws.Cells.Item(19, 4).Value = "'3.078.54"
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(20, 4).Value = "'16.35"
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(21, 4).Value = "'483.92"
$ws.Cells.Item(21, 4).ClearFormats()
$ws.Cells.Item(22, 4).Value = "'7.70"
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(23, 4).Value = "'0.684"
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Cells.Item(24, 4).Value = "'82.55"
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(25, 4).Value = "'12.79"
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(26, 4).Value = "'2.21"
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(27, 4).Value = "'10.17"
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(29, 4).Value = "'7.69"
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(30, 4).Value = "'2.30"
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Cells.Item(32, 4).Value = "'27.70"
$ws.Cells.Item(32, 4).ClearFormats()
$ws.Cells.Item(34, 4).Value = "'0.0₃0901"
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(35, 4).Value = "'1.00"
$ws.Cells.Item(35, 4).ClearFormats()
$ws.Cells.Item(36, 4).Value = "'5.67"
$ws.Cells.Item(36, 4).ClearFormats()
$ws.Cells.Item(37, 4).Value = "'0.952"
$ws.Cells.Item(37, 4).ClearFormats()
$ws.Cells.Item(38, 4).Value = "'46.24"
$ws.Cells.Item(38, 4).ClearFormats()
$ws.Cells.Item(40, 4).Value = "'1.97"
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(41, 4).Value = "'0.299"
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(42, 4).Value = "'8.28"
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(43, 4).Value = "'2.761.64"
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(44, 4).Value = "'370.27"
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(45, 4).Value = "'0.0344"
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(46, 4).Value = "'135.45"
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(47, 4).Value = "'2.49"
$ws.Cells.Item(47, 4).ClearFormats()
$ws.Cells.Item(49, 4).Value = "'24.28"
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(51, 4).Value = "'2.13"
$ws.Cells.Item(51, 4).ClearFormats()
